$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Forecast Comparison")
$ws2 = $wb.Worksheets.Item("Summary")

# --- Sheet "Forecast Comparison" updates ---
$ws1.Range("L2").Value = 1.05
$ws1.Range("L3").Value = 1.03
$ws1.Range("L4").Value = 0.82

$ws1.Range("H5").Value = 15.21
$ws1.Range("L5").Value = 0.8100000000000001

$ws1.Range("H6").Value = 14.21
$ws1.Range("L6").Value = 1.19

$ws1.Range("H7").Value = 13.23
$ws1.Range("L7").Value = 0.96

$ws1.Range("H8").Value = 12.21
$ws1.Range("L8").Value = 0.86

$ws1.Range("H9").Value = 11.21
$ws1.Range("L9").Value = 1.08

$ws1.Range("H10").Value = 10.23
$ws1.Range("L10").Value = 1.01

$ws1.Range("H11").Value = 9.23
$ws1.Range("L11").Value = 0.86

$ws1.Range("D12").Value = 10
$ws1.Range("H12").Value = 8.359999999999999
$ws1.Range("L12").Value = 0.96

$ws1.Range("H13").Value = 7.24
$ws1.Range("L13").Value = 0.97

$ws1.Range("H14").Value = 6.35
$ws1.Range("L14").Value = 1.05

$ws1.Range("H15").Value = 5.39
$ws1.Range("L15").Value = 1.1

$ws1.Range("H16").Value = 4.39
$ws1.Range("L16").Value = 1.18

$ws1.Range("H17").Value = 3.47
$ws1.Range("L17").Value = 1.16

# --- Sheet "Summary" updates ---
$ws2.Range("B10").Value = "88"
$ws2.Range("B11").Value = "43"
$ws2.Range("B14").Value = "11"
